$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: new customer record (Chat id, First name, Last name, Phone number).
# Chat id and Phone number look numeric (and the phone number even starts
# with "+"), so they must be forced to Text before assignment, otherwise
# Excel's normal General-format coercion would turn them into numbers and
# the leading "+" on the phone number would be lost.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "616525392"
$ws.Range("A22").ClearFormats()

$ws.Range("B22").Value = "Nurbek"

$ws.Range("C22").Value = "Boboyev"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "+998946696195"
$ws.Range("D22").ClearFormats()
